$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix casing of the existing "how much will spend this month..." (en_US) string
$ws.Range("D42").Value = "How much will spend this month..."

# Fill in the new LC_LOADING locale row (row 43)
$ws.Range("A43").Value = "LC_LOADING"
$ws.Range("B43").Value = "載入中..."
$ws.Range("C43").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("GOOGLETRANSLATE(B43,""ZH-TW"",""ZH-CN"")"),"载入中...")'
$ws.Range("D43").Value = "Loading..."

# New ja_JP cell E43 uses a new font/style (Arial, theme color 1)
$ws.Range("E43").Value = "読み込み中..."
$ws.Range("E43").Font.Name = "Arial"
$ws.Range("E43").Font.ThemeColor = 1
$ws.Range("E43").ReadingOrder = 0
